$d = $word.ActiveDocument

$replacements = @(
    @("323÷4=", "387÷3="),
    @("997÷4=", "829÷9="),
    @("850÷7=", "729÷7="),
    @("833÷8=", "563÷4="),
    @("385÷4=", "906÷6="),
    @("851÷5=", "359÷6="),
    @("395÷3=", "293÷4="),
    @("648÷6=", "237÷2="),
    @("209÷6=", "415÷6="),
    @("468÷9=", "433÷7="),
    @("548÷2=", "503÷6="),
    @("739÷9=", "335÷3="),
    @("149÷5=", "290÷2="),
    @("643÷6=", "533÷9="),
    @("524÷9=", "741÷3="),
    @("349÷6=", "600÷8="),
    @("740÷7=", "267÷8="),
    @("371÷3=", "755÷8="),
    @("551÷5=", "613÷5="),
    @("629÷7=", "530÷3="),
    @("292÷5=", "521÷5="),
    @("424÷8=", "118÷5="),
    @("819÷4=", "856÷2="),
    @("175÷7=", "449÷2="),
    @("565÷9=", "113÷9=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Write-Host "Done"
